$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "33.651.52"
Set-TextValue $ws.Range("E2") "  -0.87%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.770.29"
Set-TextValue $ws.Range("E3") "  -0.85%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.03%  "

# Row 5
Set-TextValue $ws.Range("D5") "223.74"
Set-TextValue $ws.Range("E5") "  +0.94%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.545"
Set-TextValue $ws.Range("E6") "  -0.82%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.00%  "

# Row 8
Set-TextValue $ws.Range("D8") "31.75"
Set-TextValue $ws.Range("E8") "  +1.03%  "

# Row 9
Set-TextValue $ws.Range("E9") "  +1.54%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0685"
Set-TextValue $ws.Range("E10") "  -4.02%  "

# Row 11
Set-TextValue $ws.Range("E11") "  +1.36%  "

# Row 12
Set-TextValue $ws.Range("D12") "2.024.46"
Set-TextValue $ws.Range("E12") "  -0.87%  "

# Row 13
Set-TextValue $ws.Range("E13") "  +4.05%  "

# Row 14
Set-TextValue $ws.Range("D14") "1.766.16"
Set-TextValue $ws.Range("E14") "  -1.22%  "

# Row 15
Set-TextValue $ws.Range("D15") "33.674.46"
Set-TextValue $ws.Range("E15") "  -0.81%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.608"
Set-TextValue $ws.Range("E16") "  -3.18%  "

# Row 17
Set-TextValue $ws.Range("E17") "  -2.50%  "

# Row 18
Set-TextValue $ws.Range("D18") "66.41"
Set-TextValue $ws.Range("E18") "  -2.30%  "

# Row 19
Set-TextValue $ws.Range("E19") "  -1.10%  "

# Row 20
Set-TextValue $ws.Range("D20") "237.89"
Set-TextValue $ws.Range("E20") "  -2.86%  "

# Row 21
Set-TextValue $ws.Range("E21") "  +0.05%  "

# Row 22
Set-TextValue $ws.Range("D22") "10.53"
Set-TextValue $ws.Range("E22") "  -1.88%  "

# Row 23
Set-TextValue $ws.Range("D23") "4.00"
Set-TextValue $ws.Range("E23") "  -1.99%  "

# Row 24
Set-TextValue $ws.Range("E24") "  -2.46%  "

# Row 25
Set-TextValue $ws.Range("D25") "159.35"
Set-TextValue $ws.Range("E25") "  +0.82%  "

# Row 26
Set-TextValue $ws.Range("D26") "16.06"
Set-TextValue $ws.Range("E26") "  -1.92%  "

# Row 27
Set-TextValue $ws.Range("E27") "  -0.12%  "

# Row 28
Set-TextValue $ws.Range("D28") "0.112"
Set-TextValue $ws.Range("E28") "  -0.24%  "

# Row 29
Set-TextValue $ws.Range("E29") "  +0.13%  "

# Row 30
Set-TextValue $ws.Range("E30") "  +1.56%  "

# Row 31
Set-TextValue $ws.Range("E31") "  -1.74%  "

# Row 32
Set-TextValue $ws.Range("E32") "  -2.71%  "

# Row 33
Set-TextValue $ws.Range("E33") "  -0.41%  "

# Row 34
Set-TextValue $ws.Range("E34") "  -1.34%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.379.86"
Set-TextValue $ws.Range("E35") "  -1.98%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.643"
Set-TextValue $ws.Range("E36") "  +0.30%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -2.40%  "

# Row 38
Set-TextValue $ws.Range("E38") "  -1.44%  "

# Row 39
Set-TextValue $ws.Range("E39") "  +5.36%  "

# Row 40
Set-TextValue $ws.Range("E40") "  +0.78%  "

# Row 41
Set-TextValue $ws.Range("D41") "77.76"
Set-TextValue $ws.Range("E41") "  -2.26%  "

# Row 42
Set-TextValue $ws.Range("B42") "ARBITRUM"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D42") "0.903"
Set-TextValue $ws.Range("E42") "  -3.62%  "

# Row 43
Set-TextValue $ws.Range("B43") "MXToken"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D43") "2.66"
Set-TextValue $ws.Range("E43") "  -2.13%  "

# Row 44
Set-TextValue $ws.Range("E44") "  +13.92%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.07"
Set-TextValue $ws.Range("E45") "  +3.84%  "

# Row 47
Set-TextValue $ws.Range("E47") "  +0.87%  "

# Row 48
Set-TextValue $ws.Range("D48") "106.87"
Set-TextValue $ws.Range("E48") "  +1.38%  "

# Row 49
Set-TextValue $ws.Range("D49") "5.81"
Set-TextValue $ws.Range("E49") "  -2.12%  "

# Row 50
Set-TextValue $ws.Range("D50") "1.924.86"

# Row 51
Set-TextValue $ws.Range("E51") "  +0.20%  "
